# Added AccountTest to see if the deposit method works.
# Updated the Test Case excel to add new case (TC004 - deposit test).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet4")

# New row 5 values (mirrors the layout of rows 2-4: TestSuiteID, TestCaseID,
# Summary, Related Requirement, Prerequisites, Test Procedure, Test Data,
# Expected Result, Actual Result, Status, Remarks, Created By, Date of
# Creation, Executed By, Date of Execution, Environment)
$ws.Range("A5").Value = "TS01"
$ws.Range("B5").Value = "TC004"
$ws.Range("C5").Value = "Verify that funds can be desposted into an account."
$ws.Range("D5").Value = "N/A"
$ws.Range("E5").Value = "1.User has created an account, either chequings or savings."
$ws.Range("F5").Value = "1. Create an account. Either chequings or savings.                                                 2. Run the deposit() method of the account and insert an amount                  3. Make sure the amount is equal to what should have been deposited."
$ws.Range("G5").Value = "1. A new Savings account object is created.                   2. A depost of 50 is made into the account."
$ws.Range("H5").Value = "1. If the depost is successful, a message letting you know that you have deposited the amount should be seen.                                                           2. If you enter a number less or equal to zero, an error message will be displayed."
$ws.Range("I5").Value = "1. If the depost is made successfully, the balance should be equal to 50 and no error message displayed.                                             2. If thedeposit in unsuccessful, the user is notified of the mistake in balance."
$ws.Range("J5").Value = "Pass"
$ws.Range("K5").Value = "Create Accounts Test Case"
$ws.Range("L5").Value = "Tyler Serio"

# M5/O5 are date-looking text ("03/24/2015", "04/13/2015") that must stay as
# plain text (shared strings), not get auto-converted into real date
# serials. Force Text format before entry, then re-apply the row's normal
# (wrap-text) style by copying formats from an existing cell already in
# that style, so no stray number-format style gets left behind.
$ws.Range("M5").NumberFormat = "@"
$ws.Range("M5").Value = "03/24/2015"
$ws.Range("M3").Copy()
$ws.Range("M5").PasteSpecial(-4122)

$ws.Range("N5").Value = "Tyler Serio"

$ws.Range("O5").NumberFormat = "@"
$ws.Range("O5").Value = "04/13/2015"
$ws.Range("O3").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("P5").Value = "OS: Windows 8.1                   IDE: Eclipse"

# Match the formatting of the other data rows: wrap text + row height 90,
# leaving the "bare" columns (A, B, D) with default style like rows 2-4.
# (M5/O5 already picked up the wrap-text style via PasteSpecial above.)
$ws.Range("C5").WrapText = $true
$ws.Range("E5:L5").WrapText = $true
$ws.Range("N5:P5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 90

# Update the view state saved in the sheet: the active cell moved to I5 and
# the visible top-left cell scrolled back to A4.
$ws.Range("I5").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
